# Correction to deal with nuclear in CES
# Change the dispatch priority of several non-nuclear electricity sources
# (hydro, onshore wind, solar PV, solar thermal, biomass, geothermal,
# offshore wind) from priority 2 down to priority 1, while leaving nuclear
# (row 4) and the other rows untouched. The dependent formulas across the
# year columns (C:AK) reference column B and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDPbES")

$ws.Range("B5").Value  = 1   # hydro
$ws.Range("B6").Value  = 1   # onshore wind
$ws.Range("B7").Value  = 1   # solar PV
$ws.Range("B8").Value  = 1   # solar thermal
$ws.Range("B9").Value  = 1   # biomass
$ws.Range("B10").Value = 1   # geothermal
$ws.Range("B14").Value = 1   # offshore wind

# Make the BDPbES sheet the active tab / selected sheet (it was "About"
# before), and leave cell B15 selected on it.
$ws.Activate()
$ws.Range("B15").Select()
